$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.094868421554565
$ws.Range("B1").Value = 2.630671739578247
$ws.Range("C1").Value = 2.745026111602783
$ws.Range("D1").Value = 3.094021797180176
$ws.Range("E1").Value = 0.7571657299995422
